$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on target cells so numeric-looking / percent-looking strings
# are stored as text (matching the original inlineStr cells) instead of being
# auto-converted to numbers by Excel.
$targetCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)
foreach ($cellAddr in $targetCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = "304.94"
$ws.Range("E2").Value = "1.48%"
$ws.Range("D3").Value = "36.21"
$ws.Range("E3").Value = "-4.93%"
$ws.Range("D4").Value = "5.031"
$ws.Range("E4").Value = "0.90%"
$ws.Range("D5").Value = "0.07816"
$ws.Range("E5").Value = "1.24%"
$ws.Range("D6").Value = "2.122"
$ws.Range("E6").Value = "-3.41%"
$ws.Range("D7").Value = "7.915"
$ws.Range("E7").Value = "-0.53%"
$ws.Range("D8").Value = "4.093"
$ws.Range("E8").Value = "2.44%"
$ws.Range("D9").Value = "0.9186"
$ws.Range("E9").Value = "0.39%"
$ws.Range("D10").Value = "0.09664"
$ws.Range("E10").Value = "6.22%"
$ws.Range("D11").Value = "0.1873"
$ws.Range("E11").Value = "4.43%"
$ws.Range("D12").Value = "0.08685"
$ws.Range("E12").Value = "2.35%"
$ws.Range("D13").Value = "0.03506"
$ws.Range("E13").Value = "-0.84%"
$ws.Range("D14").Value = "0.09909"
$ws.Range("E14").Value = "-0.26%"
$ws.Range("D15").Value = "0.001427"
$ws.Range("E15").Value = "-3.89%"
$ws.Range("D16").Value = "0.005685"
$ws.Range("E16").Value = "0.41%"
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").Value = "-0.41%"
$ws.Range("D18").Value = "2.383"
$ws.Range("E18").Value = "7.17%"
$ws.Range("D19").Value = "0.3416"
$ws.Range("E19").Value = "-1.35%"
$ws.Range("D20").Value = "0.1346"
$ws.Range("E20").Value = "2.18%"
$ws.Range("D21").Value = "4.782"
$ws.Range("E21").Value = "5.13%"
$ws.Range("D22").Value = "0.2292"
$ws.Range("E22").Value = "2.55%"
$ws.Range("D23").Value = "0.04618"
$ws.Range("E23").Value = "-0.93%"
$ws.Range("E24").Value = "15.02%"
$ws.Range("D25").Value = "0.001231"
$ws.Range("E25").Value = "0.11%"
$ws.Range("D26").Value = "0.0001402"
$ws.Range("E26").Value = "7.76%"
$ws.Range("D27").Value = "0.0004754"
$ws.Range("E27").Value = "-0.07%"
$ws.Range("D39").Value = "0.01828"
$ws.Range("E39").Value = "5.43%"
$ws.Range("D40").Value = "0.04745"
$ws.Range("E40").Value = "1.29%"
$ws.Range("D41").Value = "0.007536"
$ws.Range("E41").Value = "-3.84%"
$ws.Range("D42").Value = "0.1398"
$ws.Range("E42").Value = "0.74%"
$ws.Range("D43").Value = "0.007730"
$ws.Range("E43").Value = "0.67%"
$ws.Range("D44").Value = "0.002232"
$ws.Range("E44").Value = "-2.99%"
$ws.Range("D45").Value = "0.01026"
$ws.Range("E45").Value = "4.88%"
$ws.Range("D46").Value = "0.00006250"
$ws.Range("E46").Value = "3.60%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("D48").Value = "0.0005805"
$ws.Range("E48").Value = "0.08%"
$ws.Range("E49").Value = "184.57%"
$ws.Range("D50").Value = "0.002002"
$ws.Range("E50").Value = "-25.83%"
$ws.Range("D51").Value = "0.00002102"
$ws.Range("E51").Value = "-0.07%"

# Restore default (Normal) style so no extraneous style index is left on the cells
foreach ($cellAddr in $targetCells) {
    $ws.Range($cellAddr).Style = "Normal"
}
